$d = $word.ActiveDocument
$apos = [string][char]8217

# ---------------------------------------------------------------------------
# Helper: merge paragraph $idx1 (text becomes $newText) with the following
# paragraph $idx2, which is deleted entirely (its text + its paragraph mark).
# Processing must happen from the bottom of the document upwards so that
# paragraph indices of not-yet-processed (earlier) paragraphs stay valid.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 1) FINAL row: "A __/64" / "T __/78" (same paragraph, separated by <w:br/>)
#    -> "A __/46" / "T __/38"
# ---------------------------------------------------------------------------
$p212 = $d.Paragraphs.Item(212)
$rng212 = $p212.Range
$rng212.Find.Execute("64", $true, $false, $false, $false, $false, $true, 1, $false, "46", 2) | Out-Null
$p212b = $d.Paragraphs.Item(212)
$rng212b = $p212b.Range
$rng212b.Find.Execute("78", $true, $false, $false, $false, $false, $true, 1, $false, "38", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) SUB TOTAL row (Submission Guidelines section): "A __/16" + "T __/20"
#    (two paragraphs) -> single paragraph "__ / 16"
# ---------------------------------------------------------------------------
$p182 = $d.Paragraphs.Item(182)
$rng182 = $d.Range($p182.Range.Start, $p182.Range.End - 1)
$rng182.Text = "__ / 16"
$p183 = $d.Paragraphs.Item(183)
$p183.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Audio Presentation row TOTAL col: "A __ / 4" + "T __ / 8"
#    (two paragraphs) -> single paragraph "__ / 4"
# ---------------------------------------------------------------------------
$p175 = $d.Paragraphs.Item(175)
$rng175 = $d.Range($p175.Range.Start, $p175.Range.End - 1)
$rng175.Text = "__ / 4"
$p176 = $d.Paragraphs.Item(176)
$p176.Range.Delete()

# ---------------------------------------------------------------------------
# 4) Audio Presentation row MULTI col: "A x1" + "T x2"
#    (two paragraphs) -> single paragraph "-"
#    (replacing with a lone "-" loses run formatting in this engine, so we
#    write "-X" first and then delete the trailing "X" as a workaround)
# ---------------------------------------------------------------------------
$p173 = $d.Paragraphs.Item(173)
$rng173 = $d.Range($p173.Range.Start, $p173.Range.End - 1)
$rng173.Text = "-X"
$delRng173 = $d.Range($rng173.Start + 1, $rng173.Start + 2)
$delRng173.Delete()
$p174 = $d.Paragraphs.Item(174)
$p174.Range.Delete()

# ---------------------------------------------------------------------------
# 5) Statement 3 description paragraph - reword text
#    NOTE: Find.Execute's replacement text silently "smart-quotes" a plain
#    straight apostrophe into a curly one in this engine, but the diff wants
#    a straight apostrophe, so that particular edit is done via a direct
#    Range.Text assignment (which does NOT get auto-corrected) instead of
#    Find/Replace. We also scope the range precisely so the untouched
#    "Feel free to " / "compare and contrast" runs (incl. their proofErr
#    gramStart/gramEnd wrapper) are left alone, just like the diff.
# ---------------------------------------------------------------------------
$p129 = $d.Paragraphs.Item(129)
$p129Start = $p129.Range.Start
$p129End = $p129.Range.End - 1
$p129FullText = $d.Range($p129Start, $p129End).Text
$p129Offset = $p129FullText.IndexOf(" between your work or others.")
$p129Sub = $d.Range($p129Start + $p129Offset, $p129End)
$p129Sub.Text = " your work or others. What worked in them and what didn'tREPLACED_TAIL"
# fix up the apostrophe-adjacent tail text (avoids the engine's smart-quote
# conversion that happens specifically inside Find.Execute's replacement)
$p129Fresh = $d.Paragraphs.Item(129)
$p129FreshText = $p129Fresh.Range.Text
$p129TailOffset = $p129FreshText.IndexOf("REPLACED_TAIL")
$p129TailStart = $p129Fresh.Range.Start + $p129TailOffset
$p129TailEnd = $p129TailStart + ("REPLACED_TAIL").Length
$p129TailRng = $d.Range($p129TailStart, $p129TailEnd)
$p129TailRng.Text = "t work."

$rngAll2 = $d.Content
$rngAll2.Find.Execute("tutorials that students generated and identify some things you could", $true, $false, $false, $false, $false, $true, 1, $false, "tutorials students generated and identify what you could", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) "Statement 1" / "Statement 2" / "Statement 3" header cells - merge the
#    "Statement" run and the " N" run into a single run (text unchanged).
# ---------------------------------------------------------------------------
$p126 = $d.Paragraphs.Item(126)
$p126.Range.Find.Execute("Statement 3", $true, $false, $false, $false, $false, $true, 1, $false, "Statement 3", 2) | Out-Null

$p100 = $d.Paragraphs.Item(100)
$p100.Range.Find.Execute("Statement 2", $true, $false, $false, $false, $false, $true, 1, $false, "Statement 2", 2) | Out-Null

$p76 = $d.Paragraphs.Item(76)
$p76.Range.Find.Execute("Statement 1", $true, $false, $false, $false, $false, $true, 1, $false, "Statement 1", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Analysis, Synthesis & Evaluation SUB TOTAL: "36" -> "20", "18" -> "12"
# ---------------------------------------------------------------------------
$p73 = $d.Paragraphs.Item(73)
$p73.Range.Find.Execute("36", $true, $false, $false, $false, $false, $true, 1, $false, "20", 2) | Out-Null

$p74 = $d.Paragraphs.Item(74)
$p74.Range.Find.Execute("18", $true, $false, $false, $false, $false, $true, 1, $false, "12", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) "A poster that answers all three statements  for a general audience" -
#    merge "statements " and " for" runs into a single run (text unchanged).
# ---------------------------------------------------------------------------
$p54 = $d.Paragraphs.Item(54)
$p54.Range.Find.Execute("statements  for", $true, $false, $false, $false, $false, $true, 1, $false, "statements  for", 2) | Out-Null
